$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: reorder items inside the Walmart list cells (G10, H10, I10) ---
$ws.Range("G10").Value = '[''Mexico'', ''Wal-Mart de Mexico SAB de CV'', ''WALMART INC'', ''Walmart\xa0Inc'', ''Wal-Mart Stores'', ''Walmart Inc.'', ''Wal-Mart Stores Inc'']'
$ws.Range("H10").Value = '[''Wal-Mart de Mexico SAB de CV'', ''Walmart'', ''WALMART INC'', ''Walmart Inc'', ''Walmart Inc.'', ''Wal-Mart Stores'', ''Wal-Mart Stores Inc'', ''Walmart, Inc.'']'
$ws.Range("I10").Value = '[''Fanø'', ''Rødovre'', ''Guldborgsund'', ''Næstved'', ''Odense'', ''Randers'', ''Region\xa0Nordjylland'', ''Viborg'', ''Mariagerfjord'', ''Skive'', ''Svendborg'', ''Fredericia'', ''Greve'', ''Herning'', ''Hillerød'', ''Høje Taastrup'', ''Hørsholm'', ''Kolding'', ''Lemvig'', ''Varde'', ''Aabenraa'', ''Egedal'', ''Favrskov'', ''Region Sjælland'', ''Assens'', ''Billund'', ''Bornholm'', ''Frederiksberg'', ''Gribskov'', ''Holstebro'', ''Horsens'', ''Hvidovre'', ''Jammerbugt'', ''Lejre'', ''Lyngby-Taarbæk'', ''Nordfyn'', ''Nyborg'', ''Rebild'', ''Roskilde'', ''Skanderborg'', ''Solrød'', ''Thisted'', ''Tønder'', ''Vallensbæk'', ''Norddjurs'', ''Aalborg'', ''Faaborg-Midtfyn'', ''Hedensted'', ''Helsingør'', ''Middelfart'', ''Slagelse'', ''Syddjurs'', ''Ikast-Brande'', ''Esbjerg'', ''Aarhus'']'

# --- Row 66: reorder items inside the Daikin Industries municipality list (I66) ---
$ws.Range("I66").Value = '[''Fredericia'', ''Greve'', ''Herning'', ''Hillerød'', ''Høje Taastrup'', ''Hørsholm'', ''Kolding'', ''Lemvig'', ''Næstved'', ''Randers'', ''Skive'', ''Varde'', ''Aabenraa'', ''Furesø'', ''Ringkøbing-Skjern'', ''Horsens'', ''Mariagerfjord'', ''Skanderborg'', ''Syddjurs'', ''København'', ''Odense'', ''Esbjerg'', ''Aarhus'']'

# --- Remove row 92 (S&T Holdings Co. Ltd.) which shifts rows 93-99 up by one ---
$ws.Rows.Item(92).Delete()

Write-Output "done"
